$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1183.9727
$ws.Range("I15").Value = 1183.9727
$ws.Range("K15").Value = 3551.9181
$ws.Range("M15").Value = -3382.9181

$ws.Range("H64").Value = 3300
$ws.Range("I64").Value = 3300
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3300
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -3052
$ws.Range("N64").Value = ""

$ws.Range("H67").Value = 3300
$ws.Range("I67").Value = 3300
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3300
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -2442
$ws.Range("N67").Value = ""

$ws.Range("H137").Value = 3516.8708
$ws.Range("I137").Value = 3223.125
$ws.Range("J137").Value = 4524
$ws.Range("K137").Value = 9669.375
$ws.Range("L137").Value = 13572
$ws.Range("M137").Value = -7119.375
$ws.Range("N137").Value = -18672

$ws.Range("H138").Value = 199772.22
$ws.Range("J138").Value = 287952.16
$ws.Range("L138").Value = 863856.48
$ws.Range("N138").Value = -874136.48

$ws.Range("H141").Value = 5914.1177
$ws.Range("I141").Value = 3355.7144
$ws.Range("J141").Value = 7705
$ws.Range("K141").Value = 10067.1432
$ws.Range("L141").Value = 23115
$ws.Range("M141").Value = -4887.143199999999
$ws.Range("N141").Value = -33475


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2440.4827
$ws.Range("I22").Value = 2513.3572
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 2513.3572
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -2340.3572
$ws.Range("N22").Value = -746

$ws.Range("H99").Value = 2090
$ws.Range("I99").Value = 2184
$ws.Range("J99").Value = 1933.3334
$ws.Range("K99").Value = 2184
$ws.Range("L99").Value = 1933.3334
$ws.Range("M99").Value = -686
$ws.Range("N99").Value = -4929.3334


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 535.9
$ws.Range("I16").Value = 486.76923
$ws.Range("J16").Value = 627.1429000000001
$ws.Range("K16").Value = 486.76923
$ws.Range("L16").Value = 627.1429000000001
$ws.Range("M16").Value = -199.76923
$ws.Range("N16").Value = -1201.1429

$ws.Range("H58").Value = 1400.8235
$ws.Range("I58").Value = 1026
$ws.Range("J58").Value = 1557
$ws.Range("K58").Value = 1026
$ws.Range("L58").Value = 1557
$ws.Range("M58").Value = -823
$ws.Range("N58").Value = -1963

$ws.Range("H113").Value = 535.9
$ws.Range("I113").Value = 486.76923
$ws.Range("J113").Value = 627.1429000000001
$ws.Range("K113").Value = 486.76923
$ws.Range("L113").Value = 627.1429000000001
$ws.Range("M113").Value = 1683.23077
$ws.Range("N113").Value = -4967.1429

$ws.Range("H132").Value = 8335053
$ws.Range("I132").Value = 989.8333
$ws.Range("J132").Value = 20836148
$ws.Range("K132").Value = 2969.4999
$ws.Range("L132").Value = 62508444
$ws.Range("M132").Value = -439.4998999999998
$ws.Range("N132").Value = -62513504

$ws.Range("H134").Value = 2389
$ws.Range("I134").Value = 2301.7144
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 6905.1432
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -4370.1432
$ws.Range("N134").Value = -14070

$ws.Range("H136").Value = 1400.8235
$ws.Range("I136").Value = 1026
$ws.Range("J136").Value = 1557
$ws.Range("K136").Value = 3078
$ws.Range("L136").Value = 4671
$ws.Range("M136").Value = -528
$ws.Range("N136").Value = -9771


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 900
$ws.Range("J16").Value = 1500
$ws.Range("L16").Value = 4500
$ws.Range("N16").Value = -4846

$ws.Range("H21").Value = 4078.8572
$ws.Range("J21").Value = 4799.4
$ws.Range("L21").Value = 14398.2
$ws.Range("N21").Value = -14744.2

$ws.Range("H22").Value = 984.0741
$ws.Range("I22").Value = 958
$ws.Range("J22").Value = 990
$ws.Range("K22").Value = 2874
$ws.Range("L22").Value = 2970
$ws.Range("M22").Value = -2705
$ws.Range("N22").Value = -3308

$ws.Range("H24").Value = 3161.5386
$ws.Range("I24").Value = 1011.1111
$ws.Range("J24").Value = 8000
$ws.Range("K24").Value = 3033.3333
$ws.Range("L24").Value = 24000
$ws.Range("M24").Value = -2803.3333
$ws.Range("N24").Value = -24460

$ws.Range("H25").Value = 458.86667
$ws.Range("I25").Value = 301
$ws.Range("K25").Value = 903
$ws.Range("M25").Value = -734

$ws.Range("H26").Value = 363
$ws.Range("I26").Value = 32.857143
$ws.Range("J26").Value = 468.04544
$ws.Range("K26").Value = 98.57142899999999
$ws.Range("L26").Value = 1404.13632
$ws.Range("M26").Value = 189.428571
$ws.Range("N26").Value = -1980.13632

$ws.Range("H27").Value = 984.0741
$ws.Range("I27").Value = 958
$ws.Range("J27").Value = 990
$ws.Range("K27").Value = 2874
$ws.Range("L27").Value = 2970
$ws.Range("M27").Value = -2772
$ws.Range("N27").Value = -3174

$ws.Range("H29").Value = 111111620
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 111111620
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 333334860
$ws.Range("N29").Value = -333335414
$ws.Range("M29").Value = ""

$ws.Range("H30").Value = 458.86667
$ws.Range("I30").Value = 301
$ws.Range("K30").Value = 903
$ws.Range("M30").Value = -801

$ws.Range("H32").Value = 7754698
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 7754698
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 23264094
$ws.Range("N32").Value = -23264660
$ws.Range("M32").Value = ""

$ws.Range("H34").Value = 11111520
$ws.Range("J34").Value = 13158349
$ws.Range("L34").Value = 39475047
$ws.Range("N34").Value = -39475215


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1300.3125
$ws.Range("I107").Value = 1131.5
$ws.Range("J107").Value = 1469.125
$ws.Range("K107").Value = 1131.5
$ws.Range("L107").Value = 1469.125
$ws.Range("M107").Value = 788.5
$ws.Range("N107").Value = -5309.125

$ws.Range("H116").Value = 39999
$ws.Range("J116").Value = 39999
$ws.Range("L116").Value = 39999
$ws.Range("N116").Value = -49177

